# Agenda workbook update:
#  - Row 6 (10:00 - 10:45) used to be the "Refactoring related portions of
#    Clean Code talk" slot presented by Cory. It now becomes the
#    "Code smells and Anti-patterns presentation" slot presented by Patrick.
#  - Row 9 (01:00 - 01:30) used to be the "Presentation: ?Specific
#    refactoring techniques?" slot presented by "Patrick?". It now becomes
#    the "Refactoring related portions of Clean Code talk" slot presented
#    by Cory (i.e. that talk moved from the afternoon to the morning).
#  - The old, half-finished placeholder strings
#    ("Presentation: ?Specific refactoring techniques?" and "Patrick?")
#    are no longer used anywhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C6").Value = "Patrick"
$ws.Range("B6").Value = "Code smells and Anti-patterns presentation"

$ws.Range("B9").Value = "Refactoring related portions of Clean Code talk"
$ws.Range("C9").Value = "Cory"

$ws.Activate()
$ws.Range("C17").Select()
